$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header
$ws.Range("F1").Value = "sd_temp"

# Add new column values (sd_temp)
$values = @(
    3.27009582271967,
    2.8634823406784,
    1.77958786479798,
    3.17279898618366,
    2.58404747250405,
    1.71972595726248,
    2.68129891427354,
    2.29525515159356,
    3.00881012307577,
    2.41592976873155
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
